$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated production totals in row 3 (Lager)
$ws.Range("B3").Value = 44
$ws.Range("C3").Value = 98
$ws.Range("D3").Value = 56

# Rows 6 (Huette), 7 (Bestienjaeger) and 9 (Wollfarm) no longer carry any numbers
$ws.Range("B6:C6").ClearContents()
$ws.Range("B7:C7").ClearContents()
$ws.Range("B9:D9").ClearContents()

# Row 8 (Klantotem) also loses its numbers; clear the content but explicitly
# drop the border formatting on those cells too, so they remain present
# (blank) rather than disappearing from the sheet entirely.
$ws.Range("B8:D8").ClearContents()
$ws.Range("B8:D8").Borders.LineStyle = -4142

# Reflect the new working selection left behind by the edit
[void]$ws.Range("B7:D9").Select()
